# Auto-generated edit script: updates coin Price/Volume(1h) data and fixes a
# few reordered Coin/Link rows, matching the upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.425.59'
$ws.Range("E2").Value = '  +3.48%  '

$ws.Range("D3").Value = '2.015.07'
$ws.Range("E3").Value = '  +7.64%  '

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'0.7583"
$ws.Range("E5").Value = '  +61.23%  '

$ws.Range("D6").Value = "'257.05"
$ws.Range("E6").Value = '  +5.49%  '

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = "'0.3593"
$ws.Range("E8").Value = '  +25.10%  '

$ws.Range("D9").Value = "'29.40"
$ws.Range("E9").Value = '  +35.38%  '

$ws.Range("D10").Value = "'0.07053"
$ws.Range("E10").Value = '  +8.67%  '

$ws.Range("D11").Value = "'0.8598"
$ws.Range("E11").Value = '  +20.08%  '

$ws.Range("D12").Value = "'0.08164"

$ws.Range("D13").Value = "'102.71"
$ws.Range("E13").Value = '  +6.95%  '

$ws.Range("D14").Value = '2.017.55'
$ws.Range("E14").Value = '  +7.83%  '

$ws.Range("D15").Value = "'5.617"
$ws.Range("E15").Value = '  +9.49%  '

$ws.Range("D16").Value = "'274.15"
$ws.Range("E16").Value = '  -2.93%  '

$ws.Range("D17").Value = '31.455.69'
$ws.Range("E17").Value = '  +3.71%  '

$ws.Range("D18").Value = "'14.77"
$ws.Range("E18").Value = '  +14.03%  '

$ws.Range("D19").Value = "'5.899"
$ws.Range("E19").Value = '  +12.13%  '

$ws.Range("D20").Value = "'0.000008032"
$ws.Range("E20").Value = '  +7.07%  '

$ws.Range("D21").Value = '2.284.73'
$ws.Range("E21").Value = '  +8.11%  '

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").Value = "'0.9991"
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").Value = "'7.223"
$ws.Range("E24").Value = '  +15.63%  '

$ws.Range("D25").Value = "'10.04"
$ws.Range("E25").Value = '  +11.24%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'164.75"
$ws.Range("E26").Value = '  +0.80%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = "'0.1454"
$ws.Range("E27").Value = '  +52.04%  '

$ws.Range("D28").Value = "'20.21"
$ws.Range("E28").Value = '  +7.80%  '

$ws.Range("D29").Value = "'2.384"
$ws.Range("E29").Value = '  +26.72%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'1.595"
$ws.Range("E30").Value = '  +7.63%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'4.654"
$ws.Range("E31").Value = '  +10.23%  '

$ws.Range("D32").Value = "'1.359"
$ws.Range("E32").Value = '  +3.23%  '

$ws.Range("D33").Value = "'4.413"
$ws.Range("E33").Value = '  +7.25%  '

$ws.Range("D34").Value = "'0.05247"
$ws.Range("E34").Value = '  +8.89%  '

$ws.Range("D35").Value = "'1.243"
$ws.Range("E35").Value = '  +11.00%  '

$ws.Range("D36").Value = "'0.7692"
$ws.Range("E36").Value = '  +12.01%  '

$ws.Range("D37").Value = "'2.783"
$ws.Range("E37").Value = '  +2.56%  '

$ws.Range("D38").Value = "'0.02030"
$ws.Range("E38").Value = '  +7.45%  '

$ws.Range("D39").Value = "'2.954"
$ws.Range("E39").Value = '  +5.20%  '

$ws.Range("D40").Value = "'6.809"
$ws.Range("E40").Value = '  +9.72%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = "'80.37"
$ws.Range("E41").Value = '  +6.99%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'2.206"
$ws.Range("E42").Value = '  +15.24%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = "'0.4771"
$ws.Range("E43").Value = '  +13.79%  '

$ws.Range("D44").Value = "'0.8641"
$ws.Range("E44").Value = '  +4.85%  '

$ws.Range("D45").Value = "'105.22"
$ws.Range("E45").Value = '  +4.62%  '

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = '  +0.22%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = "'7.771"
$ws.Range("E47").Value = '  +11.20%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'10.07"
$ws.Range("E48").Value = '  +4.12%  '

$ws.Range("D49").Value = "'0.4406"
$ws.Range("E49").Value = '  +12.92%  '

$ws.Range("D50").Value = "'37.14"
$ws.Range("E50").Value = '  +5.80%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = "'948.59"
$ws.Range("E51").Value = '  +5.45%  '
